# Add a new column (K) of data for year 2020, mirroring the existing
# year columns (D:J) that already hold 2013-2019 data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 (the thin rule row above the header) --------------------------
# K2 should look exactly like J2 (blank, bottom border, vertically centered).
$ws.Range("J2").Copy() | Out-Null
$ws.Range("K2").PasteSpecial(-4122) | Out-Null   # xlPasteFormats

# --- Row 3 (header row with the year numbers) -----------------------------
$ws.Range("J3").Copy() | Out-Null
$ws.Range("K3").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$ws.Range("K3").Value = 2020

# --- Row 4 ------------------------------------------------------------------
$ws.Range("J4").Copy() | Out-Null
$ws.Range("K4").PasteSpecial(-4122) | Out-Null
$ws.Range("K4").Value = 0

# --- Row 5 -------------------------------------------------------------------
$ws.Range("I5").Copy() | Out-Null
$ws.Range("K5").PasteSpecial(-4122) | Out-Null
$ws.Range("K5").Value = 48.2

# --- Row 6 ---------------------------------------------------------------
$ws.Range("J6").Copy() | Out-Null
$ws.Range("K6").PasteSpecial(-4122) | Out-Null
$ws.Range("K6").Value = 19.3

# --- Row 7 ---------------------------------------------------------------
$ws.Range("J7").Copy() | Out-Null
$ws.Range("K7").PasteSpecial(-4122) | Out-Null
$ws.Range("K7").Value = 24.2

# --- Row 8 ---------------------------------------------------------------
$ws.Range("J8").Copy() | Out-Null
$ws.Range("K8").PasteSpecial(-4122) | Out-Null
$ws.Range("K8").Value = 8.3

$excel.CutCopyMode = 0

# --- Row 1 header height changed slightly once the new column was added --
$ws.Rows.Item(1).RowHeight = 63.75

# --- Selection left behind by the editor (purely cosmetic) ---------------
$ws.Range("J22").Select() | Out-Null
